$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Cadm1", "Cadm3", "ECs", 3, 1, 4.716522666666667, 14.149568, 0.530211572117814, 0.530211572117814, 3, 1, 4.370833, 13.112499, 0.0996525630698175, 0.09965256306981748, 20.61513291671467, 185.536196250432, 0.05283694213081755, 0.05283694213081754),
  @("ECs", "Cadm1", "Cadm3", "FAPs", 3, 1, 4.716522666666667, 14.149568, 0.530211572117814, 0.530211572117814, 3, 1, 37.55834333333333, 112.67503, 0.8563093528905953, 0.8563093528905952, 177.1447776541156, 1594.30299888704, 0.4540251282153105, 0.4540251282153104),
  @("ECs", "Cadm1", "Cadm3", "sCs", 3, 1, 4.716522666666667, 14.149568, 0.530211572117814, 0.530211572117814, 3, 1, 1.931542, 5.794626, 0.04403808403958729, 0.04403808403958728, 9.110161624618668, 81.991454621568, 0.02334950177168599, 0.02334950177168599),
  @("FAPs", "Cadm1", "Cadm3", "ECs", 2, 0.6666666666666666, 0.236474, 0.709422, 0.02658340904223816, 0.02658340904223817, 3, 1, 4.370833, 13.112499, 0.0996525630698175, 0.09965256306981748, 1.033588362842, 9.302295265578, 0.002649104846192395, 0.002649104846192395),
  @("FAPs", "Cadm1", "Cadm3", "FAPs", 2, 0.6666666666666666, 0.236474, 0.709422, 0.02658340904223816, 0.02658340904223817, 3, 1, 37.55834333333333, 112.67503, 0.8563093528905953, 0.8563093528905952, 8.881571681406665, 79.93414513265999, 0.02276362179458496, 0.02276362179458497),
  @("FAPs", "Cadm1", "Cadm3", "sCs", 2, 0.6666666666666666, 0.236474, 0.709422, 0.02658340904223816, 0.02658340904223817, 3, 1, 1.931542, 5.794626, 0.04403808403958729, 0.04403808403958728, 0.456759462908, 4.110835166172, 0.001170682401460809, 0.001170682401460809),
  @("sCs", "Cadm1", "Cadm3", "ECs", 3, 1, 3.942551666666667, 11.827655, 0.4432050188399478, 0.4432050188399479, 3, 1, 4.370833, 13.112499, 0.0996525630698175, 0.09965256306981748, 17.23223492887167, 155.090114359845, 0.04416651609280755, 0.04416651609280756),
  @("sCs", "Cadm1", "Cadm3", "FAPs", 3, 1, 3.942551666666667, 11.827655, 0.4432050188399478, 0.4432050188399479, 3, 1, 37.55834333333333, 112.67503, 0.8563093528905953, 0.8563093528905952, 148.0757091060722, 1332.68138195465, 0.3795206028806998, 0.3795206028806998),
  @("sCs", "Cadm1", "Cadm3", "sCs", 3, 1, 3.942551666666667, 11.827655, 0.4432050188399478, 0.4432050188399479, 3, 1, 1.931542, 5.794626, 0.04403808403958729, 0.04403808403958728, 7.615204131336667, 68.53683718203, 0.01951789986644049, 0.01951789986644049)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
